# D4_monthly_prevalence_vaccination codebook - "Data Model" sheet
#
# The commit removes the standalone "dose" row (it was a leftover single
# VarName entry with no Description/Format/Vocabulary) from the vaccine
# coverage codebook, and documents the Rule for the last variable
# (PP_month) by adding its formula in the "Rule" (column K) cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)   # "Data Model" sheet

# Row 7 holds the orphaned "dose" VarName (A7="dose", nothing else useful).
# Deleting the whole row shifts rows 8:15 up to 7:14, which also drops
# "dose" out of the shared-string table automatically.
$ws.Rows.Item(7).Delete()

# The last row is now row 14 (PP_month). Document its Rule in column K.
$ws.Cells.Item(14, 11).Value = "(Vacc_observed_before_month + Vacc_observed_month)/NFUP_month"

# Keep the view/selection state sane after the row deletion (cursor was
# sitting on row 2, move it down one row to A3 to mirror the saved file).
$ws.Activate()
$ws.Range("A3").Select()
